$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format from the previous date cell so the new date cells
# reuse the existing "date" cell style instead of creating a new one.
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 8: new timesheet entry
$ws.Range("A8").Value = 43861
$ws.Range("B8").Value = 2.4
$ws.Range("C8").Value = "Access USPS data"

# Row 9: new timesheet entry
$ws.Range("A9").Value = 43862
$ws.Range("B9").Value = 1.9
$ws.Range("C9").Value = "Access USPS data"

# Update selection to mimic post-edit cursor position
$ws.Range("A10").Select()
